# Weekly update: a new price observation is published for this market/product.
# The new record is inserted as row 17 (pushing the existing rows 17-230 down
# to 18-231), carrying forward the same Volumen/Precio/Unidad/Origen values as
# the entry that is now directly below it (row 18, the old row 17), but with
# a brand-new "Fecha" (2021-12-22, Excel serial 44552).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..230 down to 18..231, creating a blank row 17.
$ws.Rows.Item(17).Insert()

# Seed the new row 17 with the same data as the row right after it (the data
# that used to live in row 17 before the shift), then overwrite its date.
$srcRow = $ws.Range("A18:R18")
$newRow = $ws.Range("A17:R17")
$newRow.Value = $srcRow.Value()

# New "Fecha" for the newly inserted observation: 2021-12-22 (serial 44552).
$ws.Cells.Item(17, 4).Value = 44552
